# Make email fields lowercase on submit — append the 4 new test rows that
# exercise the lower-casing behaviour (rows 33-36 on Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the last existing data row (32) as a formatting template for the four
# new rows so the new cells inherit the same styles (bold/bordered ID column,
# date-formatted "Date Sent" column, etc.) as the rest of the table.
$ws.Range("A32:I32").Copy($ws.Range("A33:I33"))
$ws.Range("A32:I32").Copy($ws.Range("A34:I34"))
$ws.Range("A32:I32").Copy($ws.Range("A35:I35"))
$ws.Range("A32:I32").Copy($ws.Range("A36:I36"))

# Row 33
$ws.Range("A33").Value = 31
$ws.Range("B33").Value = "Hi Luke"
$ws.Range("C33").Value = 45009.48006944444
$ws.Range("D33").Value = "This is my first email as the developer."
$ws.Range("E33").Value = "chelly-xox@hotmail.co.uk"
$ws.Range("F33").Value = "Micah Chuku"
$ws.Range("G33").Value = "lukevaughn@aol.com"
$ws.Range("H33").Value = "Luke Vaughn"
$ws.Range("I33").Value = $false

# Row 34
$ws.Range("A34").Value = 32
$ws.Range("B34").Value = "Hi again, Luke"
$ws.Range("C34").Value = 45009.48297453704
# "'Sup bro?" genuinely starts with an apostrophe character in the source
# data (not a formatting quote-prefix). A plain `.Value =` assignment is
# interpreted the same way Excel interprets interactive typing, which
# strips a leading "'" as a text-qualifier. Route it through
# Copy/PasteSpecial(values) instead so the literal apostrophe is kept as
# real cell content, matching the inline string in the target workbook.
$scratch = $ws.Range("Z100")
$scratch.Formula = "=""'Sup bro?"""
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$scratch.Clear()
$excel.CutCopyMode = $false
$ws.Range("E34").Value = "chelly-xox@hotmail.co.uk"
$ws.Range("F34").Value = "Micah Chuku"
$ws.Range("G34").Value = "lukevaughn@aol.com"
$ws.Range("H34").Value = "Luke Vaughn"
$ws.Range("I34").Value = $true

# Row 35
$ws.Range("A35").Value = 33
$ws.Range("B35").Value = "Hello Luke"
$ws.Range("C35").Value = 45009.48892361111
$ws.Range("D35").Value = "I am writing this email with your email address capitalised. Hopefully the web app should turn this to lowercase before submitting to the API. Have a good day!"
$ws.Range("E35").Value = "n.schneider@gmail.com"
$ws.Range("F35").Value = "Natalia Schneider"
$ws.Range("G35").Value = "lukevaughn@aol.com"
$ws.Range("H35").Value = "Luke Vaughn"
$ws.Range("I35").Value = $false

# Row 36
$ws.Range("A36").Value = 34
$ws.Range("B36").Value = "Hi again again"
$ws.Range("C36").Value = 45009.49577546296
$ws.Range("D36").Value = "Testing out after removing patter on email field"
$ws.Range("E36").Value = "chelly-xox@hotmail.co.uk"
$ws.Range("F36").Value = "Micah Chuku"
$ws.Range("G36").Value = "lukevaughn@aol.com"
$ws.Range("H36").Value = "Luke Vaughn"
$ws.Range("I36").Value = $false
